$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/ada-tooth"
$ws.Range("B3").Value = "8.0.0"
$ws.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$ws.Range("B9").Value = "LinuxForHealth Team"
